$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "59.381.95"
$ws.Range("E2").Value = "  +2.79%  "

Set-TextValue $ws "D3" "2.593.17"
$ws.Range("E3").Value = "  +1.92%  "

$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue $ws "D5" "520.69"
$ws.Range("E5").Value = "  +0.48%  "

Set-TextValue $ws "D6" "139.32"
$ws.Range("E6").Value = "  +0.63%  "

Set-TextValue $ws "D7" "0.999"
$ws.Range("E7").Value = "  -0.19%  "

Set-TextValue $ws "D8" "0.567"
$ws.Range("E8").Value = "  +0.88%  "

Set-TextValue $ws "D9" "2.611.28"
$ws.Range("E9").Value = "  +2.69%  "

Set-TextValue $ws "D10" "6.50"
$ws.Range("E10").Value = "  -0.22%  "

$ws.Range("E11").Value = "  +2.27%  "

Set-TextValue $ws "D12" "0.332"
$ws.Range("E12").Value = "  +2.38%  "

$ws.Range("E13").Value = "  +1.74%  "

Set-TextValue $ws "D14" "3.057.53"
$ws.Range("E14").Value = "  +2.00%  "

Set-TextValue $ws "D15" "59.277.81"
$ws.Range("E15").Value = "  +2.64%  "

Set-TextValue $ws "D16" "20.38"
$ws.Range("E16").Value = "  +2.10%  "

Set-TextValue $ws "D17" "2.618.78"
$ws.Range("E17").Value = "  +3.15%  "

Set-TextValue $ws "D18" "0.0000133"
$ws.Range("E18").Value = "  +0.31%  "

Set-TextValue $ws "D19" "339.28"
$ws.Range("E19").Value = "  +1.91%  "

Set-TextValue $ws "D20" "4.32"
$ws.Range("E20").Value = "  +1.16%  "

Set-TextValue $ws "D21" "10.18"
$ws.Range("E21").Value = "  +0.83%  "

Set-TextValue $ws "D22" "6.49"
$ws.Range("E22").Value = "  +6.18%  "

Set-TextValue $ws "D23" "0.998"
$ws.Range("E23").Value = "  -0.21%  "

Set-TextValue $ws "D24" "66.52"
$ws.Range("E24").Value = "  +2.45%  "

$ws.Range("E25").Value = "  +1.73%  "

Set-TextValue $ws "D26" "0.404"
$ws.Range("E26").Value = "  +0.86%  "

Set-TextValue $ws "D27" "0.997"
$ws.Range("E27").Value = "  -0.24%  "

$ws.Range("E28").Value = "  +1.77%  "

$ws.Range("E29").Value = "  -0.14%  "

Set-TextValue $ws "D30" "0.0₃0724"
$ws.Range("E30").Value = "  -3.46%  "

Set-TextValue $ws "D31" "5.96"
$ws.Range("E31").Value = "  -3.50%  "

Set-TextValue $ws "D32" "18.82"
$ws.Range("E32").Value = "  +2.05%  "

Set-TextValue $ws "D33" "1.57"
$ws.Range("E33").Value = "  +0.14%  "

Set-TextValue $ws "D34" "149.31"
$ws.Range("E34").Value = "  +0.27%  "

Set-TextValue $ws "D35" "3.99"
$ws.Range("E35").Value = "  +0.82%  "

Set-TextValue $ws "D36" "1.13"
$ws.Range("E36").Value = "  +0.29%  "

Set-TextValue $ws "D37" "36.27"
$ws.Range("E37").Value = "  +1.62%  "

$ws.Range("E38").Value = "  +3.83%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws "D39" "0.830"
$ws.Range("E39").Value = "  +1.24%  "

$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws "D40" "0.824"
$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("E41").Value = "  +2.32%  "

Set-TextValue $ws "D42" "0.997"
$ws.Range("E42").Value = "  -0.24%  "

Set-TextValue $ws "D43" "275.06"
$ws.Range("E43").Value = "  +5.75%  "

Set-TextValue $ws "D44" "10.74"
$ws.Range("E44").Value = "  +0.89%  "

$ws.Range("E45").Value = "  +2.57%  "

Set-TextValue $ws "D46" "0.0951"
$ws.Range("E46").Value = "  -0.14%  "

Set-TextValue $ws "D47" "0.0521"
$ws.Range("E47").Value = "  +0.35%  "

Set-TextValue $ws "D48" "18.48"
$ws.Range("E48").Value = "  +0.21%  "

Set-TextValue $ws "D49" "1.984.17"
$ws.Range("E49").Value = "  +0.34%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D50" "0.0220"
$ws.Range("E50").Value = "  -0.45%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D51" "4.47"
$ws.Range("E51").Value = "  -0.57%  "
